$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update the raw input data cells (columns E/F/G) that feed the
#    formulas further down the sheet.
# ---------------------------------------------------------------
$ws.Range("E21").Value = 88214000

$ws.Range("F25").Value = 20000

$ws.Range("F27").Value = 2699000

$ws.Range("F28").Value = 2518000
$ws.Range("G28").Value = 1000

$ws.Range("F29").Value = 76000

$ws.Range("F31").Value = 76000

$ws.Range("F35").Value = 2523000

$ws.Range("F36").Value = 2528000

$ws.Range("F37").Value = 1000

$ws.Range("F40").Value = 2762000

$ws.Range("F41").Value = 47000

$ws.Range("F66").Value = 5871000

$ws.Range("F67").Value = 76000

$ws.Range("F68").Value = 7523000

# ---------------------------------------------------------------
# 2) Fill in the previously-empty N column formulas for the ovens
#    (rows 66-68), mirroring the existing M-column formulas.
# ---------------------------------------------------------------
$ws.Range("N66").Formula = "=G66*H66/(J66*H66*480*L66)"
$ws.Range("N67").Formula = "=G67*H67/(J67*H67*480*L67)"
$ws.Range("N68").Formula = "=G68*H68/(J68*H68*480*L68)"

# ---------------------------------------------------------------
# 3) Rebuild the small "Gem. ploegen" summary block (M19:N22).
#    First add the new header row (19) re-using the existing
#    "Gem. ploegen na persen" text, then clear the old labels in
#    row 70 so their shared-string slots are freed up and reused
#    in-place by the new "ovens" labels, matching how the workbook
#    was edited originally.
# ---------------------------------------------------------------
$ws.Range("M19").Value = "Gem. ploegen na persen"
$ws.Range("M19").Font.Bold = $true

$ws.Range("L70").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("L70").Font.Bold = $false
$ws.Range("M70").Font.Bold = $false
$ws.Range("L70").Font.ThemeColor = 1
$ws.Range("M70").Font.ThemeColor = 1

$ws.Range("M21").Value = "Gem. ploegen na ovens"
$ws.Range("N21").Value = "Gem. ploegen voor ovens"
$ws.Range("N21").Font.Bold = $true

$ws.Range("N19").Value = "Gem. ploegen voor persen"
$ws.Range("N19").Font.Bold = $true

$ws.Range("M20").NumberFormat = "0.00"
$ws.Range("M20").Formula = "=AVERAGE(M25:M41)"
$ws.Range("N20").NumberFormat = "0.00"
$ws.Range("N20").Formula = "=AVERAGE(N25:N41)"

$ws.Range("M22").Formula = "=AVERAGE(M66:M68)"
$ws.Range("N22").NumberFormat = "0.00"
$ws.Range("N22").Formula = "=AVERAGE(N66:N68)"

# ---------------------------------------------------------------
# 4) Misc layout tweaks
# ---------------------------------------------------------------
$ws.Columns(14).ColumnWidth = 25

$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("N23").Select()
